# Apply the updated "cryptos" price/volume snapshot to Sheet1.
#
# For each affected row, only the columns that actually changed are listed
# (B = Coin name, C = Link, D = Price, E = Volume(1h)). Rows 50/51 also swap
# the ImmutableX/Filecoin entries (B/C/D all change, not just the numbers).
#
# Columns D and E hold free-form text in the source data (e.g. "98.192.77",
# "3.368.13" -- thousands-grouped values using a literal dot separator -- and
# "  +4.33%  " -- percentages padded with spaces). Excel auto-converts a
# plain decimal-looking string (like "254.98") into a real number when
# assigned via .Value, which would change the stored cell type away from the
# original text. To keep these as text, numeric-looking values are entered
# with a leading apostrophe (forces text entry) and the cell style is reset
# to "Normal" right after so no stray number-format/quote-prefix styling is
# left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = [ordered]@{
    2 = @{ "D"="98.192.77"; "E"="  +4.33%  " }
    3 = @{ "D"="3.368.13"; "E"="  +9.63%  " }
    4 = @{ "E"="  -0.05%  " }
    5 = @{ "D"="254.98"; "E"="  +8.03%  " }
    6 = @{ "D"="623.40"; "E"="  +2.24%  " }
    7 = @{ "E"="  +8.38%  " }
    8 = @{ "D"="0.386"; "E"="  +2.33%  " }
    9 = @{ "E"="  -0.03%  " }
    10 = @{ "D"="3.365.66"; "E"="  +9.51%  " }
    11 = @{ "D"="0.813"; "E"="  +0.37%  " }
    12 = @{ "E"="  +1.35%  " }
    13 = @{ "D"="97.996.00"; "E"="  +4.22%  " }
    14 = @{ "D"="35.89"; "E"="  +5.90%  " }
    15 = @{ "E"="  +2.83%  " }
    16 = @{ "D"="3.991.66"; "E"="  +9.24%  " }
    17 = @{ "E"="  +3.47%  " }
    18 = @{ "D"="3.368.23"; "E"="  +9.74%  " }
    19 = @{ "D"="3.67"; "E"="  +3.35%  " }
    20 = @{ "D"="14.81"; "E"="  +3.14%  " }
    21 = @{ "D"="481.88"; "E"="  +8.58%  " }
    22 = @{ "D"="5.90"; "E"="  +2.99%  " }
    23 = @{ "E"="  +10.36%  " }
    24 = @{ "D"="9.23"; "E"="  +4.59%  " }
    25 = @{ "D"="5.72"; "E"="  +3.88%  " }
    26 = @{ "D"="88.17"; "E"="  +4.08%  " }
    27 = @{ "D"="12.04"; "E"="  +0.90%  " }
    28 = @{ "D"="3.542.63"; "E"="  +9.25%  " }
    29 = @{ "E"="  -0.15%  " }
    30 = @{ "E"="  +5.04%  " }
    31 = @{ "E"="  +0.88%  " }
    32 = @{ "D"="0.125"; "E"="  +2.98%  " }
    33 = @{ "D"="1.00"; "E"="  -0.02%  " }
    34 = @{ "D"="9.27"; "E"="  +3.73%  " }
    35 = @{ "D"="27.37"; "E"="  +7.60%  " }
    36 = @{ "D"="525.41"; "E"="  +9.04%  " }
    37 = @{ "E"="  +1.37%  " }
    38 = @{ "D"="7.35"; "E"="  -2.97%  " }
    39 = @{ "D"="1.95"; "E"="  +3.67%  " }
    40 = @{ "D"="24.81"; "E"="  +3.19%  " }
    41 = @{ "E"="  +3.09%  " }
    42 = @{ "D"="3.82"; "E"="  +2.68%  " }
    43 = @{ "E"="  +1.74%  " }
    44 = @{ "E"="  +5.56%  " }
    45 = @{ "D"="0.793"; "E"="  +18.16%  " }
    46 = @{ "E"="  -0.01%  " }
    47 = @{ "E"="  -0.24%  " }
    48 = @{ "E"="  +6.60%  " }
    49 = @{ "D"="45.54"; "E"="  +4.29%  " }
    50 = @{ "B"="Filecoin"; "C"="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; "D"="4.55"; "E"="  +7.04%  " }
    51 = @{ "B"="ImmutableX"; "C"="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; "D"="1.37"; "E"="  +6.64%  " }
}

$numericPattern = '^-?[0-9]*\.?[0-9]+$'

foreach ($row in $rowUpdates.Keys) {
    $cols = $rowUpdates[$row]
    foreach ($col in $cols.Keys) {
        $ref = "$col$row"
        $text = $cols[$col]
        $needsTextGuard = ($col -eq "D") -and ($text -match $numericPattern)
        if ($needsTextGuard) {
            $ws.Range($ref).Value = "'" + $text
            $ws.Range($ref).Style = "Normal"
        } else {
            $ws.Range($ref).Value = $text
        }
    }
}

Write-Output "Updated $($rowUpdates.Count) rows in $($ws.Name)"
